$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("processed_new_data")
$rng = $ws1.Range("A1")
$rng.Interior.ThemeColor = 10
